$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 229 (shifts existing rows 229..249 down to 230..250)
$ws.Rows.Item(229).Insert()

# Populate the newly inserted row 229 with the new price record
$ws.Cells.Item(229, 1).Value = 7
$ws.Cells.Item(229, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(229, 3).Value = "Ñuble"
$ws.Cells.Item(229, 4).Value = 44826
$ws.Cells.Item(229, 5).Value = 16
$ws.Cells.Item(229, 6).Value = "Fruta"
$ws.Cells.Item(229, 7).Value = 100108
$ws.Cells.Item(229, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(229, 9).Value = 100108005
$ws.Cells.Item(229, 10).Value = "Piña"
$ws.Cells.Item(229, 11).Value = "Caramelo"
$ws.Cells.Item(229, 12).Value = "Segunda"
$ws.Cells.Item(229, 13).Value = 120
$ws.Cells.Item(229, 14).Value = 21000
$ws.Cells.Item(229, 15).Value = 22000
$ws.Cells.Item(229, 16).Value = 21500
$ws.Cells.Item(229, 17).Value = "`$/caja 14 unidades"
$ws.Cells.Item(229, 18).Value = "Ecuador"
$ws.Cells.Item(229, 19).Value = 1536
$ws.Cells.Item(229, 20).Value = 14

# Make sure the date cell keeps the date number format used by the rest of column D
$ws.Cells.Item(229, 4).NumberFormat = $ws.Cells.Item(230, 4).NumberFormat
